$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "15 TL - 15 TL"
$ws.Range("G3").Value = "39,87 TRY - 79,76 TRY - 797,68 TRY"
$ws.Range("G4").Value = "27,84 TRY - 55,69 TRY - 398,83 TRY"
$ws.Range("G5").Value = "7,97 TRY - 15,96 TRY - 199,41 TRY"
$ws.Range("G6").Value = "8.300,01 TL - 99,71 TL"
$ws.Range("G7").Value = "1 TRY (Kredi kartı ile ödemelerde ek olarak nakit avans faizi uygulanır.)"
$ws.Range("G8").Value = "19,94 TRY - 39,88 TRY - 398,84 TRY"
$ws.Range("G9").Value = "13,92 TRY - 27,85 TRY - 199,42 TRY"
$ws.Range("G10").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("G11").Value = "3,99 TRY - 7,98 TRY - 99,71 TRY"
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"
$ws.Range("F13").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = "8.300 TL - 7,97 TL"
